# Append a new data row (row 53) to Sheet1, mirroring the existing
# "date / weekday / hour / ranking" log table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 53

# Column A holds the date as literal text (e.g. "2025/10/02") in every
# existing row, not a real Excel date serial. Force text interpretation
# via a temporary "@" number format so Excel doesn't auto-convert the
# "2025/10/03" string into a date value, then drop the formatting again
# so the new cell matches the unstyled cells around it.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "2025/10/03"
$ws.Range("A" + $newRow).ClearFormats()

$ws.Range("B" + $newRow).Value = "金"
$ws.Range("C" + $newRow).Value = 2
$ws.Range("D" + $newRow).Value = 3
